$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly sampling was recorded. It becomes the new row 10, pushing the
# existing rows 10-26 down to 11-27 (each of those keeps its original values,
# just one row lower), and the former row 26 lands on the new row 27.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the latest data point.
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Macroferia Regional de Talca"
$ws.Range("C10").Value = "Maule"
$ws.Range("D10").Value = 44467
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 100112026
$ws.Range("G10").Value = "Haba"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 9000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 9000
$ws.Range("N10").Value = '$/saco 25 kilos'
$ws.Range("O10").Value = "Región de O'Higgins"
$ws.Range("P10").Value = 360
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
